# KHL injuries snapshot refresh (2025-12-21 04:0x UTC scrape / 12:xx +08:00 local).
#
# 1) "snapshot": every still-injured player keeps their row but gets a fresh
#    `scraped_at` (column K) timestamp from the new scrape pass.
# 2) Two players who were on "snapshot" (CSKA's Buchelnikov & Moiseev) have
#    recovered, so their rows are removed from "snapshot" and appended (with
#    a RETURN status + changed_at/changed_day) to the "returned" sheet,
#    replacing the now-stale rows that previously lived there.
# 3) Two newly-injured Kunlun Red Star players (Bischoff & Grolo) fill the
#    vacated rows at the bottom of "snapshot".

$wb = $excel.ActiveWorkbook

$snapshot = $wb.Worksheets.Item("snapshot")
$returned = $wb.Worksheets.Item("returned")

# --- 1) Refresh scraped_at (column K) for the players who are still injured ---
$scrapedUpdates = @(
    @{ Row = 2;  Value = "2025-12-21T04:08:40.074273+00:00" },
    @{ Row = 3;  Value = "2025-12-21T04:08:42.480320+00:00" },
    @{ Row = 4;  Value = "2025-12-21T04:08:42.480361+00:00" },
    @{ Row = 5;  Value = "2025-12-21T04:08:42.480383+00:00" },
    @{ Row = 6;  Value = "2025-12-21T04:08:44.959073+00:00" },
    @{ Row = 7;  Value = "2025-12-21T04:08:47.454990+00:00" },
    @{ Row = 8;  Value = "2025-12-21T04:08:49.907935+00:00" },
    @{ Row = 9;  Value = "2025-12-21T04:08:49.907963+00:00" },
    @{ Row = 10; Value = "2025-12-21T04:08:52.290441+00:00" },
    @{ Row = 11; Value = "2025-12-21T04:08:56.650714+00:00" },
    @{ Row = 12; Value = "2025-12-21T04:08:59.183833+00:00" },
    @{ Row = 13; Value = "2025-12-21T04:09:01.219794+00:00" },
    @{ Row = 14; Value = "2025-12-21T04:09:06.181926+00:00" },
    @{ Row = 15; Value = "2025-12-21T04:09:06.181945+00:00" },
    @{ Row = 16; Value = "2025-12-21T04:09:06.181953+00:00" },
    @{ Row = 17; Value = "2025-12-21T04:09:06.181961+00:00" },
    @{ Row = 18; Value = "2025-12-21T04:09:08.723582+00:00" },
    @{ Row = 19; Value = "2025-12-21T04:09:08.723604+00:00" },
    @{ Row = 20; Value = "2025-12-21T04:09:08.723612+00:00" },
    @{ Row = 21; Value = "2025-12-21T04:09:10.635988+00:00" },
    @{ Row = 22; Value = "2025-12-21T04:09:10.636015+00:00" },
    @{ Row = 23; Value = "2025-12-21T04:09:10.636033+00:00" },
    @{ Row = 24; Value = "2025-12-21T04:09:13.125286+00:00" },
    @{ Row = 25; Value = "2025-12-21T04:09:13.125307+00:00" },
    @{ Row = 26; Value = "2025-12-21T04:09:15.056007+00:00" },
    @{ Row = 27; Value = "2025-12-21T04:09:15.056043+00:00" },
    @{ Row = 28; Value = "2025-12-21T04:09:15.056065+00:00" },
    @{ Row = 29; Value = "2025-12-21T04:09:16.969311+00:00" },
    @{ Row = 30; Value = "2025-12-21T04:09:19.062701+00:00" },
    @{ Row = 31; Value = "2025-12-21T04:09:19.062731+00:00" }
)

foreach ($u in $scrapedUpdates) {
    $snapshot.Cells.Item($u.Row, 11).Value = $u.Value
}

# --- 2) Two CSKA players recovered: move them off "snapshot" onto "returned" ---

# Capture their snapshot data (rows 32-33) before overwriting the rows.
$returnedPlayers = @(
    @{
        TeamAbbr  = "ЦСК"
        TeamName  = "ЦСКА"
        PlayerName = "Бучельников Дмитрий"
        PlayerUid = "1369_ЦСК_бучельниковдмитрий"
    },
    @{
        TeamAbbr  = "ЦСК"
        TeamName  = "ЦСКА"
        PlayerName = "Моисеев Данила"
        PlayerUid = "1369_ЦСК_моисеевданила"
    }
)

# --- 3) Two newly-injured Kunlun Red Star players take over rows 32-33 ---
$newInjured = @(
    @{
        TeamAbbr   = "ШДР"
        TeamName   = "Драконы"
        TeamSlug   = "kunlun"
        PlayerName = "Бишофф Джейк"
        Number     = "28"
        Position   = "защитник"
        PlayerId   = "45490"
        PlayerUid  = "1369_ШДР_бишоффджейк"
        SourceUrl  = "https://www.khl.ru/clubs/kunlun/team/"
        ScrapedAt  = "2025-12-21T04:10:07.846979+00:00"
    },
    @{
        TeamAbbr   = "ШДР"
        TeamName   = "Драконы"
        TeamSlug   = "kunlun"
        PlayerName = "Гроло Жереми"
        Number     = "75"
        Position   = "защитник"
        PlayerId   = "45343"
        PlayerUid  = "1369_ШДР_гроложереми"
        SourceUrl  = "https://www.khl.ru/clubs/kunlun/team/"
        ScrapedAt  = "2025-12-21T04:10:07.847014+00:00"
    }
)

$row = 32
foreach ($p in $newInjured) {
    $snapshot.Cells.Item($row, 1).Value  = $p.TeamAbbr
    $snapshot.Cells.Item($row, 2).Value  = $p.TeamName
    $snapshot.Cells.Item($row, 3).Value  = $p.TeamSlug
    $snapshot.Cells.Item($row, 4).Value  = $p.PlayerName
    # Number/PlayerId are numeric-looking but must stay text (like every
    # other row in this column) - quote-prefix forces text entry, then
    # resetting the style back to Normal drops the quote-prefix style flag
    # so the cell ends up with no explicit style, matching the rest of the sheet.
    $snapshot.Cells.Item($row, 5).Value  = "'" + $p.Number
    $snapshot.Cells.Item($row, 5).Style  = "Normal"
    $snapshot.Cells.Item($row, 6).Value  = $p.Position
    $snapshot.Cells.Item($row, 7).Value  = "'" + $p.PlayerId
    $snapshot.Cells.Item($row, 7).Style  = "Normal"
    $snapshot.Cells.Item($row, 8).Value  = $p.PlayerUid
    $snapshot.Cells.Item($row, 9).Value  = "injured_active"
    $snapshot.Cells.Item($row, 10).Value = $p.SourceUrl
    $snapshot.Cells.Item($row, 11).Value = $p.ScrapedAt
    $row = $row + 1
}

# The two trailing rows (34-35) that used to hold the Kunlun players are now
# redundant — their data has already been promoted into rows 32-33 above.
$snapshot.Range("A34:K35").EntireRow.Delete()

# --- Replace the stale "returned" rows with the two CSKA players who just recovered ---
$changedAt  = "2025-12-21T12:10:08.350201+08:00"
$changedDay = "2025-12-21"

$row = 2
foreach ($p in $returnedPlayers) {
    $returned.Cells.Item($row, 1).Value = $p.TeamAbbr
    $returned.Cells.Item($row, 2).Value = $p.TeamName
    $returned.Cells.Item($row, 3).Value = $p.PlayerName
    $returned.Cells.Item($row, 4).Value = $p.PlayerUid
    $returned.Cells.Item($row, 5).Value = "RETURN"
    $returned.Cells.Item($row, 6).Value = $changedAt
    # changed_day ("2025-12-21") looks like a date, so Excel would otherwise
    # auto-convert it to a date serial. Quote-prefix keeps it literal text,
    # then resetting the style to Normal drops the quote-prefix style flag.
    $returned.Cells.Item($row, 7).Value = "'" + $changedDay
    $returned.Cells.Item($row, 7).Style = "Normal"
    $row = $row + 1
}

# The previous row 4 (Sochi's Bikmullin) is dropped entirely.
$returned.Range("A4:G4").EntireRow.Delete()
